$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")
$ws.Unprotect()
$ws.Range("A1").Value = "test"
